$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 762.5
$ws.Range("I31").Value = 762.5
$ws.Range("K31").Value = 2287.5
$ws.Range("M31").Value = -2057.5
$ws.Range("H53").Value = 568.0476
$ws.Range("I53").Value = 319.16666
$ws.Range("J53").Value = 899.8889
$ws.Range("K53").Value = 319.16666
$ws.Range("L53").Value = 899.8889
$ws.Range("M53").Value = 317.83334
$ws.Range("N53").Value = -2173.8889
$ws.Range("H125").Value = 1745.3334
$ws.Range("I125").Value = 2362.3333
$ws.Range("J125").Value = 1436.8334
$ws.Range("K125").Value = 21260.9997
$ws.Range("L125").Value = 12931.5006
$ws.Range("M125").Value = -18800.9997
$ws.Range("N125").Value = -17851.5006
$ws.Range("H126").Value = 41835
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 41835
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 41835
$ws.Range("N126").Value = -51715
$ws.Range("H127").Value = 2261.6
$ws.Range("I127").Value = 351.75
$ws.Range("J127").Value = 2956.0908
$ws.Range("K127").Value = 1055.25
$ws.Range("L127").Value = 8868.2724
$ws.Range("M127").Value = 3904.75
$ws.Range("N127").Value = -18788.2724
$ws.Range("H128").Value = 41835
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 41835
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 41835
$ws.Range("N128").Value = -51795
$ws.Range("H129").Value = 834.54
$ws.Range("I129").Value = 355.85715
$ws.Range("J129").Value = 870.5699
$ws.Range("K129").Value = 1067.57145
$ws.Range("L129").Value = 2611.7097
$ws.Range("M129").Value = 3932.42855
$ws.Range("N129").Value = -12611.7097
$ws.Range("H130").Value = 42495.555
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 42495.555
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 42495.555
$ws.Range("N130").Value = -52535.555
$ws.Range("H131").Value = 2390.3635
$ws.Range("I131").Value = 498.75
$ws.Range("J131").Value = 3471.2856
$ws.Range("K131").Value = 1496.25
$ws.Range("L131").Value = 10413.8568
$ws.Range("M131").Value = 3543.75
$ws.Range("N131").Value = -20493.8568
$ws.Range("H132").Value = 28576360
$ws.Range("I132").Value = 32262756
$ws.Range("J132").Value = 6800
$ws.Range("K132").Value = 96788268
$ws.Range("L132").Value = 20400
$ws.Range("M132").Value = -96785738
$ws.Range("N132").Value = -25460
$ws.Range("H133").Value = 42393.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 42393.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 42393.332
$ws.Range("N133").Value = -52513.332
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 1048.5
$ws.Range("I135").Value = 475.07693
$ws.Range("J135").Value = 3533.3333
$ws.Range("K135").Value = 4275.69237
$ws.Range("L135").Value = 31799.9997
$ws.Range("M135").Value = -1740.69237
$ws.Range("N135").Value = -36869.9997
$ws.Range("H136").Value = 45126.668
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 45126.668
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 45126.668
$ws.Range("N136").Value = -55326.668
$ws.Range("H137").Value = 1986615.4
$ws.Range("I137").Value = 2802986.2
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 8408958.600000001
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -8406408.600000001
$ws.Range("N137").Value = -17100
$ws.Range("H138").Value = 2548.65
$ws.Range("I138").Value = 687.2105
$ws.Range("J138").Value = 2985.284
$ws.Range("K138").Value = 2061.6315
$ws.Range("L138").Value = 8955.852000000001
$ws.Range("M138").Value = 3078.3685
$ws.Range("N138").Value = -19235.852
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 47670.715
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 47670.715
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 47670.715
$ws.Range("N140").Value = -58030.715
$ws.Range("H141").Value = 24466.334
$ws.Range("I141").Value = 26899.625
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 80698.875
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -75518.875
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1958.2307
$ws.Range("I61").Value = 1400
$ws.Range("K61").Value = 1400
$ws.Range("M61").Value = -1188
$ws.Range("H74").Value = 10947.2
$ws.Range("I74").Value = 17010
$ws.Range("J74").Value = 4884.4
$ws.Range("K74").Value = 17010
$ws.Range("L74").Value = 4884.4
$ws.Range("M74").Value = -16136
$ws.Range("N74").Value = -6632.4
$ws.Range("H77").Value = 10947.2
$ws.Range("I77").Value = 17010
$ws.Range("J77").Value = 4884.4
$ws.Range("K77").Value = 85050
$ws.Range("L77").Value = 24422
$ws.Range("M77").Value = -80682
$ws.Range("N77").Value = -33158
$ws.Range("H132").Value = 1956.7407
$ws.Range("I132").Value = 1182.6818
$ws.Range("K132").Value = 3548.0454
$ws.Range("M132").Value = -1018.0454
$ws.Range("H136").Value = 1958.2307
$ws.Range("I136").Value = 1400
$ws.Range("K136").Value = 4200
$ws.Range("M136").Value = -1650

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 35016.332
$ws.Range("J82").Value = 35016.332
$ws.Range("L82").Value = 35016.332
$ws.Range("N82").Value = -35782.332
$ws.Range("H85").Value = 35016.332
$ws.Range("J85").Value = 35016.332
$ws.Range("L85").Value = 35016.332
$ws.Range("N85").Value = -37668.332
$ws.Range("H134").Value = 2433.5
$ws.Range("I134").Value = 1921.8334
$ws.Range("J134").Value = 5503.5
$ws.Range("K134").Value = 5765.5002
$ws.Range("L134").Value = 16510.5
$ws.Range("M134").Value = -3230.5002
$ws.Range("N134").Value = -21580.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3297.238
$ws.Range("I31").Value = 1080.9333
$ws.Range("J31").Value = 8838
$ws.Range("K31").Value = 1080.9333
$ws.Range("L31").Value = 8838
$ws.Range("M31").Value = -785.9332999999999
$ws.Range("N31").Value = -9428
$ws.Range("H34").Value = 3297.238
$ws.Range("I34").Value = 1080.9333
$ws.Range("J34").Value = 8838
$ws.Range("K34").Value = 1080.9333
$ws.Range("L34").Value = 8838
$ws.Range("M34").Value = -878.9332999999999
$ws.Range("N34").Value = -9242
$ws.Range("H58").Value = 2669.2395
$ws.Range("I58").Value = 1643.2142
$ws.Range("J58").Value = 6499.7334
$ws.Range("K58").Value = 1643.2142
$ws.Range("L58").Value = 6499.7334
$ws.Range("M58").Value = -1440.2142
$ws.Range("N58").Value = -6905.7334
$ws.Range("H132").Value = 3605.0527
$ws.Range("I132").Value = 2299.3
$ws.Range("J132").Value = 5055.8887
$ws.Range("K132").Value = 6897.900000000001
$ws.Range("L132").Value = 15167.6661
$ws.Range("M132").Value = -4367.900000000001
$ws.Range("N132").Value = -20227.6661
$ws.Range("H134").Value = 8223
$ws.Range("I134").Value = 10461.546
$ws.Range("J134").Value = 4119
$ws.Range("K134").Value = 31384.638
$ws.Range("L134").Value = 12357
$ws.Range("M134").Value = -28849.638
$ws.Range("N134").Value = -17427
$ws.Range("H136").Value = 2669.2395
$ws.Range("I136").Value = 1643.2142
$ws.Range("J136").Value = 6499.7334
$ws.Range("K136").Value = 4929.642599999999
$ws.Range("L136").Value = 19499.2002
$ws.Range("M136").Value = -2379.642599999999
$ws.Range("N136").Value = -24599.2002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 25905
$ws.Range("J53").Value = 25905
$ws.Range("L53").Value = 25905
$ws.Range("N53").Value = -27167
$ws.Range("H126").Value = 3024.35
$ws.Range("I126").Value = 2762.775
$ws.Range("J126").Value = 4070.65
$ws.Range("K126").Value = 8288.325000000001
$ws.Range("L126").Value = 12211.95
$ws.Range("M126").Value = -5818.325000000001
$ws.Range("N126").Value = -17151.95

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5252.8
$ws.Range("I132").Value = 2426.5
$ws.Range("J132").Value = 13025.125
$ws.Range("K132").Value = 7279.5
$ws.Range("L132").Value = 39075.375
$ws.Range("M132").Value = -4749.5
$ws.Range("N132").Value = -44135.375
$ws.Range("H136").Value = 5647.2666
$ws.Range("I136").Value = 1868.1666
$ws.Range("J136").Value = 8166.6665
$ws.Range("K136").Value = 5604.4998
$ws.Range("L136").Value = 24499.9995
$ws.Range("M136").Value = -3054.4998
$ws.Range("N136").Value = -29599.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 37396
$ws.Range("J58").Value = 50094
$ws.Range("L58").Value = 50094
$ws.Range("N58").Value = -50710
$ws.Range("H132").Value = 11116888
$ws.Range("I132").Value = 5990.227
$ws.Range("J132").Value = 41671856
$ws.Range("K132").Value = 17970.681
$ws.Range("L132").Value = 125015568
$ws.Range("M132").Value = -15440.681
$ws.Range("N132").Value = -125020628
$ws.Range("H136").Value = 17176.637
$ws.Range("J136").Value = 10122.25
$ws.Range("L136").Value = 30366.75
$ws.Range("N136").Value = -35466.75
